$wb = $excel.ActiveWorkbook

# --- Update the Stream_seepage sheet: insert a new top data row (row 2) ---
$ws3 = $wb.Worksheets.Item("Stream_seepage")

$ws3.Rows("2:2").Insert()

# Copy the date-cell formatting from the row below (now row 3) so the new
# A2 cell gets the same short-date number format as the rest of column A.
$ws3.Range("A3").Copy()
$ws3.Range("A2").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws3.Range("A2").Value = 45183
$ws3.Range("B2").Value = "HOB fit"
$ws3.Range("C2").Value = @'
With the final round of updates, we see that MW_19 over simulates drawdown in the summer while MW17 matches more closely and these wells are directly on opposite sides of the rivers. In MODFLOW there is likely enough connectivity that they see similar heads, but the difference might be developed if the GDE mapping was adjusted. Decreasing extinction depth from 10 m to 8m didn't show a noticeable impact on HOB fit for some reason. 
-> I continued by decreasing the ext_dp for woodlan to 6m (forest still at 8m) and didn't see a head fit change. The mean ET out did decline from 124k to 104k to 76k m3/day with SFR in declining with similar amounts (177k, 165k, 147k) to account for the change with like storage impacts as well. (I had been loading the wrong hob output)
-> I updated the hob path read and the model fit greatly improved (from the starting point NSE went from 0.51 to .65 and RMSE from 2.0 to 1.7 m)
'@
$ws3.Range("D2").Value = @'
Test a slight further decline in ext_dp then proceed (woodland 5 m with forest at 8 m, and riparian scrub at 3 m, slightly improved fit, best to test across all 100 next)
'@

# --- Update the active sheet / selections to match the new state ---
$ws2 = $wb.Worksheets.Item("setback_distance")
[void]$ws2.Range("D3").Select()

[void]$ws3.Activate()
[void]$ws3.Range("D3").Select()
